# Feature Tracker update
# - Reworks the feature list: adds several new requested features
#   (per-die modifiers, roll groupings, add color, roll variants),
#   and re-sorts existing rows so the table reads in the new order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Start clean: wipe out the old table contents (including the now-unused
# column F) so we can lay out the refreshed table from scratch.
$ws.Range("A1:F9").ClearContents()

# Header row
$ws.Cells.Item(1,1).Value = "Feature"
$ws.Cells.Item(1,2).Value = "Description"
$ws.Cells.Item(1,3).Value = "Completed Version"
$ws.Cells.Item(1,4).Value = "Requestors"

# Save Custom Rolls
$ws.Cells.Item(2,1).Value = "Save Custom Rolls"
$ws.Cells.Item(2,2).Value = "Type in roll in custom roll tab, hit save, appears in saved custom rolls tab, pressing on it will roll everything within."
$ws.Cells.Item(2,3).Value = "1.4.1"
$ws.Cells.Item(2,4).Value = "_dharwin - Reddit /r/dnd"
$ws.Cells.Item(2,5).Value = "Nitrogen06 - Reddit /r/rpg"

# Fate Dice
$ws.Cells.Item(3,1).Value = "Fate Dice "
$ws.Cells.Item(3,2).Value = "Have a custom type of dice that rolls between -1 and 1"
$ws.Cells.Item(3,3).Value = "1.3.0"
$ws.Cells.Item(3,4).Value = "joethomp - Reddit /r/rpg"

# Have per die modifiers (new)
$ws.Cells.Item(4,1).Value = "Have per die modifiers"
$ws.Cells.Item(4,2).Value = "Instead of having a single modifier that is on a per roll basis, have it on a per die basis"
$ws.Cells.Item(4,4).Value = "Nitrogen06 - Reddit /r/rpg"

# Roll groupings (new) - description jotted down first, feature name filled in after
$ws.Cells.Item(5,2).Value = "When you have many different saved rolls it can be overwhelming, groups of rolls would be nice"

# Advantage/Disadvantage
$ws.Cells.Item(6,1).Value = "Advantage/Disadvantage"
$ws.Cells.Item(6,2).Value = "Roll a set of dice twice and take the higher of the two sets of rolls"
$ws.Cells.Item(6,4).Value = "Kevin Ryan - silverghost2@gmail.com"
$ws.Cells.Item(6,5).Value = "UraniumKnight - Reddit /r/rpg"

# Drop X High/Low
$ws.Cells.Item(7,1).Value = "Drop X High/Low"
$ws.Cells.Item(7,2).Value = "Roll a set of dice once and drop the highest X number of rolls or lowest X number of rolls"
$ws.Cells.Item(7,4).Value = "Kevin Ryan - silverghost2@gmail.com"
$ws.Cells.Item(7,5).Value = "UraniumKnight - Reddit /r/rpg"

# Add Color (new)
$ws.Cells.Item(8,1).Value = "Add Color"
$ws.Cells.Item(8,2).Value = "Having all the dice be in gray scale is not fun to look at"
$ws.Cells.Item(8,4).Value = "Kris Fiala"

# Back to row 5 to fill in the feature name
$ws.Cells.Item(5,1).Value = "Roll groupings"

# Roll Variants (new)
$ws.Cells.Item(9,1).Value = "Roll Variants"
$ws.Cells.Item(9,2).Value = "Allow for setting variant versions of a roll, i.e. d6 vs d6(poison) vs d6(green)"
$ws.Cells.Item(9,4).Value = "Weston Fiala"

# Match the saved selection state (cursor sat on A10 after entering data)
$ws.Range("A10").Select()
